$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.288.11'
$ws.Cells.Item(2, 5).Value = '  +0.36%  '
$ws.Cells.Item(3, 4).Value = '3.572.92'
$ws.Cells.Item(3, 5).Value = '  +2.59%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).Value = '608.93'
$ws.Cells.Item(5, 5).Value = '  +0.76%  '
$ws.Cells.Item(6, 4).Value = '145.55'
$ws.Cells.Item(6, 5).Value = '  +1.75%  '
$ws.Cells.Item(7, 4).Value = '3.569.78'
$ws.Cells.Item(7, 5).Value = '  +2.54%  '
$ws.Cells.Item(8, 5).Value = '  +0.23%  '
$ws.Cells.Item(10, 5).Value = '  +1.40%  '
$ws.Cells.Item(11, 5).Value = '  -3.19%  '
$ws.Cells.Item(12, 5).Value = '  +0.65%  '
$ws.Cells.Item(13, 4).Value = '4.180.05'
$ws.Cells.Item(13, 5).Value = '  +2.61%  '
$ws.Cells.Item(14, 5).Value = '  +2.63%  '
$ws.Cells.Item(15, 4).Value = '29.99'
$ws.Cells.Item(15, 5).Value = '  -1.16%  '
$ws.Cells.Item(16, 4).Value = '3.586.25'
$ws.Cells.Item(16, 5).Value = '  +2.93%  '
$ws.Cells.Item(17, 4).Value = '66.390.97'
$ws.Cells.Item(17, 5).Value = '  +0.39%  '
$ws.Cells.Item(18, 5).Value = '  -0.96%  '
$ws.Cells.Item(19, 4).Value = '11.54'
$ws.Cells.Item(19, 5).Value = '  +11.10%  '
$ws.Cells.Item(20, 5).Value = '  +1.19%  '
$ws.Cells.Item(21, 4).Value = '14.88'
$ws.Cells.Item(21, 5).Value = '  +0.92%  '
$ws.Cells.Item(22, 4).Value = '429.91'
$ws.Cells.Item(22, 5).Value = '  +1.98%  '
$ws.Cells.Item(23, 5).Value = '  +4.64%  '
$ws.Cells.Item(24, 4).Value = '79.22'
$ws.Cells.Item(24, 5).Value = '  +2.25%  '
$ws.Cells.Item(25, 4).Value = '3.716.26'
$ws.Cells.Item(25, 5).Value = '  +2.67%  '
$ws.Cells.Item(26, 5).Value = '  +0.04%  '
$ws.Cells.Item(27, 4).Value = '0.0000119'
$ws.Cells.Item(27, 5).Value = '  +4.18%  '
$ws.Cells.Item(28, 5).Value = '  +2.43%  '
$ws.Cells.Item(29, 4).Value = '7.95'
$ws.Cells.Item(29, 5).Value = '  -0.35%  '
$ws.Cells.Item(30, 4).Value = '9.09'
$ws.Cells.Item(30, 5).Value = '  -2.54%  '
$ws.Cells.Item(31, 5).Value = '  +0.15%  '
$ws.Cells.Item(32, 4).Value = '25.66'
$ws.Cells.Item(32, 5).Value = '  +2.07%  '
$ws.Cells.Item(33, 5).Value = '  -1.28%  '
$ws.Cells.Item(34, 4).Value = '3.569.81'
$ws.Cells.Item(34, 5).Value = '  +2.58%  '
$ws.Cells.Item(35, 5).Value = '  -5.43%  '
$ws.Cells.Item(37, 5).Value = '  +1.64%  '
$ws.Cells.Item(38, 5).Value = '  +2.68%  '
$ws.Cells.Item(39, 4).Value = '5.61'
$ws.Cells.Item(39, 5).Value = '  +0.69%  '
$ws.Cells.Item(40, 4).Value = '177.51'
$ws.Cells.Item(40, 5).Value = '  +4.09%  '
$ws.Cells.Item(41, 5).Value = '  -0.02%  '
$ws.Cells.Item(42, 5).Value = '  -1.58%  '
$ws.Cells.Item(43, 5).Value = '  +2.79%  '
$ws.Cells.Item(44, 4).Value = '0.897'
$ws.Cells.Item(44, 5).Value = '  +0.94%  '
$ws.Cells.Item(45, 4).Value = '1.94'
$ws.Cells.Item(45, 5).Value = '  +1.22%  '
$ws.Cells.Item(46, 4).Value = '46.17'
$ws.Cells.Item(46, 5).Value = '  +2.35%  '
$ws.Cells.Item(47, 5).Value = '  +1.30%  '
$ws.Cells.Item(48, 4).Value = '25.71'
$ws.Cells.Item(48, 5).Value = '  -1.62%  '
$ws.Cells.Item(49, 5).Value = '  +2.84%  '
$ws.Cells.Item(50, 4).Value = '7.15'
$ws.Cells.Item(50, 5).Value = '  +0.48%  '
$ws.Cells.Item(51, 4).Value = '23.53'
$ws.Cells.Item(51, 5).Value = '  +9.41%  '